$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a 'clean' number to Excel's auto-detection
# (e.g. '317.62', '1.0000', '10.90') must be forced to stay TEXT -- exactly as
# they were stored before the edit (t="inlineStr" / shared-string, never t="n").
# A leading apostrophe is the standard Excel quote-prefix trick: it forces the
# literal text into the cell (dropping the apostrophe itself) without touching the
# cell's number format.

$ws.Range("D2").Value = "'24.970.42"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "'1.709.54"
$ws.Range("E3").Value = "  -0.33%  "
$ws.Range("E4").Value = "  -0.36%  "
$ws.Range("D5").Value = "'317.62"
$ws.Range("E5").Value = "  -0.12%  "
$ws.Range("D6").Value = "'1.0000"
$ws.Range("E6").Value = "  -0.27%  "
$ws.Range("D7").Value = "'0.4030"
$ws.Range("E7").Value = "  +1.50%  "
$ws.Range("D8").Value = "'0.4079"
$ws.Range("E8").Value = "  -0.93%  "
$ws.Range("B9").Value = "Polygon"
$ws.Range("C9").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D9").Value = "'1.482"
$ws.Range("E9").Value = "  -2.99%  "
$ws.Range("B10").Value = "OKB"
$ws.Range("C10").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D10").Value = "'53.85"
$ws.Range("E10").Value = "  +0.82%  "
$ws.Range("E11").Value = "  -0.36%  "
$ws.Range("D12").Value = "'0.08836"
$ws.Range("E12").Value = "  -1.14%  "
$ws.Range("D13").Value = "'26.41"
$ws.Range("E13").Value = "  +5.66%  "
$ws.Range("D14").Value = "'7.518"
$ws.Range("E14").Value = "  -2.83%  "
$ws.Range("D15").Value = "'8.145"
$ws.Range("E15").Value = "  -0.14%  "
$ws.Range("D16").Value = "'0.00001362"
$ws.Range("E16").Value = "  -0.64%  "
$ws.Range("D17").Value = "'1.749.40"
$ws.Range("E17").Value = "  +4.19%  "
$ws.Range("D18").Value = "'97.01"
$ws.Range("E18").Value = "  -3.55%  "
$ws.Range("D19").Value = "'0.07175"
$ws.Range("E19").Value = "  +0.30%  "
$ws.Range("D20").Value = "'21.15"
$ws.Range("E20").Value = "  +4.88%  "
$ws.Range("D21").Value = "'7.272"
$ws.Range("E21").Value = "  -3.20%  "
$ws.Range("D22").Value = "'1.0000"
$ws.Range("E22").Value = "  -0.38%  "
$ws.Range("D23").Value = "'14.42"
$ws.Range("E23").Value = "  -0.78%  "
$ws.Range("D24").Value = "'24.971.20"
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").Value = "'2.330"
$ws.Range("E25").Value = "  -0.59%  "
$ws.Range("D26").Value = "'2.915"
$ws.Range("E26").Value = "  -7.17%  "
$ws.Range("D27").Value = "'23.32"
$ws.Range("E27").Value = "  +0.25%  "
$ws.Range("D28").Value = "'6.242"
$ws.Range("E28").Value = "  +18.93%  "
$ws.Range("D29").Value = "'166.98"
$ws.Range("E29").Value = "  +1.08%  "
$ws.Range("D30").Value = "'146.84"
$ws.Range("E30").Value = "  +4.18%  "
$ws.Range("D31").Value = "'8.406"
$ws.Range("E31").Value = "  -9.03%  "
$ws.Range("D32").Value = "'1.930.44"
$ws.Range("E32").Value = "  +3.41%  "
$ws.Range("D33").Value = "'2.228"
$ws.Range("E33").Value = "  +13.24%  "
$ws.Range("D34").Value = "'0.08891"
$ws.Range("E34").Value = "  -1.45%  "
$ws.Range("D35").Value = "'0.03242"
$ws.Range("E35").Value = "  +7.82%  "
$ws.Range("D36").Value = "'7.274"
$ws.Range("E36").Value = "  -7.81%  "
$ws.Range("D37").Value = "'1.032"
$ws.Range("E37").Value = "  -5.13%  "
$ws.Range("D38").Value = "'0.2870"
$ws.Range("E38").Value = "  +2.15%  "
$ws.Range("D39").Value = "'0.8489"
$ws.Range("E39").Value = "  +4.19%  "
$ws.Range("D40").Value = "'10.90"
$ws.Range("E40").Value = "  -2.43%  "
$ws.Range("D41").Value = "'0.09355"
$ws.Range("E41").Value = "  +0.36%  "
$ws.Range("D42").Value = "'14.32"
$ws.Range("E42").Value = "  -1.90%  "
$ws.Range("D43").Value = "'1.478"
$ws.Range("E43").Value = "  -0.84%  "
$ws.Range("D44").Value = "'17.65"
$ws.Range("E44").Value = "  +5.19%  "
$ws.Range("D45").Value = "'2.741"
$ws.Range("E45").Value = "  +3.39%  "
$ws.Range("D46").Value = "'0.7467"
$ws.Range("E46").Value = "  +1.21%  "
$ws.Range("D47").Value = "'4.255"
$ws.Range("E47").Value = "  -0.43%  "
$ws.Range("D48").Value = "'1.402"
$ws.Range("E48").Value = "  +3.98%  "
$ws.Range("D49").Value = "'0.9995"
$ws.Range("E49").Value = "  -0.23%  "
$ws.Range("D50").Value = "'142.36"
$ws.Range("E50").Value = "  +0.97%  "
$ws.Range("D51").Value = "'0.08382"
$ws.Range("E51").Value = "  +3.57%  "
